$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.817.97"
$ws.Range("E2").Value = "  +5.00%  "
$ws.Range("D3").Value = "3.116.09"
$ws.Range("E3").Value = "  +3.34%  "
$ws.Range("E4").Value = "  -0.09%  "
$cell = $ws.Range("D5")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "584.11"
$cell.Style = $origStyle
$cell = $ws.Range("D6")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "144.32"
$cell.Style = $origStyle
$ws.Range("E7").Value = "  -0.10%  "
$ws.Range("D8").Value = "3.108.24"
$ws.Range("E8").Value = "  +3.40%  "
$ws.Range("E9").Value = "  +1.51%  "
$ws.Range("E10").Value = "  +10.79%  "
$cell = $ws.Range("D11")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "5.75"
$cell.Style = $origStyle
$ws.Range("E11").Value = "  +8.21%  "
$ws.Range("E12").Value = "  +1.37%  "
$ws.Range("E13").Value = "  +6.13%  "
$ws.Range("E14").Value = "  +3.95%  "
$ws.Range("D16").Value = "3.627.64"
$ws.Range("E16").Value = "  +3.06%  "
$ws.Range("E17").Value = "  -0.58%  "
$ws.Range("D18").Value = "3.111.53"
$ws.Range("E18").Value = "  +3.04%  "
$ws.Range("D19").Value = "62.765.56"
$ws.Range("E19").Value = "  +4.93%  "
$cell = $ws.Range("D20")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "465.95"
$cell.Style = $origStyle
$ws.Range("E20").Value = "  +6.15%  "
$ws.Range("E21").Value = "  +2.81%  "
$ws.Range("E22").Value = "  +0.89%  "
$cell = $ws.Range("D23")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "7.56"
$cell.Style = $origStyle
$ws.Range("E23").Value = "  +6.10%  "
$cell = $ws.Range("D24")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "13.32"
$cell.Style = $origStyle
$ws.Range("E24").Value = "  -0.31%  "
$cell = $ws.Range("D25")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "82.17"
$cell.Style = $origStyle
$ws.Range("E25").Value = "  +1.59%  "
$ws.Range("E26").Value = "  +0.00%  "
$cell = $ws.Range("D27")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "2.24"
$cell.Style = $origStyle
$ws.Range("E27").Value = "  -0.35%  "
$ws.Range("E28").Value = "  +4.76%  "
$ws.Range("E29").Value = "  -0.01%  "
$cell = $ws.Range("D30")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "8.26"
$cell.Style = $origStyle
$ws.Range("E30").Value = "  +5.45%  "
$cell = $ws.Range("D31")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "6.81"
$cell.Style = $origStyle
$ws.Range("E31").Value = "  +7.72%  "
$ws.Range("E32").Value = "  +6.58%  "
$ws.Range("E33").Value = "  +3.30%  "
$ws.Range("D34").Value = "0.0₃0854"
$ws.Range("E34").Value = "  +7.71%  "
$ws.Range("E35").Value = "  +11.62%  "
$ws.Range("E36").Value = "  +3.64%  "
$ws.Range("E37").Value = "  +1.87%  "
$cell = $ws.Range("D38")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "3.23"
$cell.Style = $origStyle
$ws.Range("E38").Value = "  +14.55%  "
$cell = $ws.Range("D39")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "50.94"
$cell.Style = $origStyle
$ws.Range("E39").Value = "  +3.51%  "
$cell = $ws.Range("D40")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "432.44"
$cell.Style = $origStyle
$ws.Range("E40").Value = "  +6.37%  "
$ws.Range("E41").Value = "  +1.25%  "
$ws.Range("D42").Value = "2.939.67"
$ws.Range("E42").Value = "  +5.86%  "
$ws.Range("E43").Value = "  +4.10%  "
$ws.Range("E44").Value = "  +9.08%  "
$cell = $ws.Range("D45")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.111"
$cell.Style = $origStyle
$ws.Range("E45").Value = "  +3.57%  "
$ws.Range("E46").Value = "  +6.37%  "
$cell = $ws.Range("D47")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "35.37"
$cell.Style = $origStyle
$ws.Range("E47").Value = "  +4.00%  "
$ws.Range("E48").Value = "  -0.04%  "
$cell = $ws.Range("D49")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "123.25"
$cell.Style = $origStyle
$ws.Range("E49").Value = "  +0.25%  "
$ws.Range("E50").Value = "  +0.32%  "
$cell = $ws.Range("D51")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "24.73"
$cell.Style = $origStyle
$ws.Range("E51").Value = "  +4.43%  "
